# Swap the match data between row 46 and row 47 (the "id" column A stays put,
# as do the columns that already hold identical values in both rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","F","G","I","J","K","L","M","N","O","P","Q","U","V","W","X","AB")

foreach ($col in $cols) {
    $cell46 = $ws.Range("$col" + "46")
    $cell47 = $ws.Range("$col" + "47")
    $val46 = $cell46.Value()
    $val47 = $cell47.Value()
    $cell46.Value = $val47
    $cell47.Value = $val46
}
